$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A (Requirement names) top to bottom
$ws.Range("A2").Value = "Reflection Paper 1"
$ws.Range("A3").Value = "Reflection Paper 2"
$ws.Range("A4").Value = "Reflection Paper 3"
$ws.Range("A5").Value = "Reflection Paper 4"
$ws.Range("A6").Value = "Participation"

# Fill column B (Points)
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 20

# Fill column C (Percent)
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 20
$ws.Range("C5").Value = 20
$ws.Range("C6").Value = 20

# Fill column D (Due Dates)
$ws.Range("D2").Value = "Friday, September 18"
$ws.Range("D2").NumberFormat = "d-mmm"
$ws.Range("D3").Value = "Friday, October 9"
$ws.Range("D4").Value = "Friday, November 13"
$ws.Range("D5").Value = "Friday, December 11"
$ws.Range("D6").Value = "Rolling Basis"

# Column widths (best-fit approximation)
$ws.Range("A1").ColumnWidth = 15.166666666666666
$ws.Range("D1").ColumnWidth = 17.709635416666668

# Final selection, matching the saved cursor position
$ws.Range("D6").Select()
